$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShareSkill")

# --- Column widths (columns 3-13 get new widths) ---
# Target raw widths: C:D=27.140625, E:H=18.5703125, I:J=22.140625, K:M=16
# The ColumnWidth setter quantizes to 1/6 character-width steps, so we pick
# the closest achievable ColumnWidth values.
$ws.Columns.Item(3).ColumnWidth = 26.333333333333332
$ws.Columns.Item(4).ColumnWidth = 26.333333333333332
$ws.Columns.Item(5).ColumnWidth = 17.666666666666668
$ws.Columns.Item(6).ColumnWidth = 17.666666666666668
$ws.Columns.Item(7).ColumnWidth = 17.666666666666668
$ws.Columns.Item(8).ColumnWidth = 17.666666666666668
$ws.Columns.Item(9).ColumnWidth = 21.333333333333332
$ws.Columns.Item(10).ColumnWidth = 21.333333333333332
$ws.Columns.Item(11).ColumnWidth = 15.166666666666666
$ws.Columns.Item(12).ColumnWidth = 15.166666666666666
$ws.Columns.Item(13).ColumnWidth = 15.166666666666666

# --- Apply the header style (same as existing header cells, fill s="1") to the
# new header cells before writing their values ---
$ws.Range("A1").Copy()
$ws.Range("C1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 1 headers (final target layout) ---
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Description"
$ws.Range("C1").Value = "CategoryId"
$ws.Range("D1").Value = "SubcategoryId"
$ws.Range("E1").Value = "Tags"
$ws.Range("F1").Value = "Service Type"
$ws.Range("G1").Value = "Location Type"
$ws.Range("H1").Value = "StartDate"
$ws.Range("I1").Value = "EndDate"
$ws.Range("J1").Value = "StartTime"
$ws.Range("K1").Value = "EndTime"
$ws.Range("L1").Value = "Skill Trade"
$ws.Range("M1").Value = "Skills"
$ws.Range("N1").Value = "Credit"
$ws.Range("O1").Value = "State"

# --- Row 2 data (final target layout) ---
$ws.Range("A2").Value = "Software Tester"
$ws.Range("B2").Value = "Can teach QA Skills"
$ws.Range("C2").Value = "Programming & Tech"
$ws.Range("D2").Value = "QA"
$ws.Range("E2").Value = "test"
$ws.Range("F2").Value = "Hourlybasis"
$ws.Range("G2").Value = "Online"
$ws.Range("H2").Value = 17082019
$ws.Range("I2").Value = 23082019
$ws.Range("J2").Value = "0100PM"
$ws.Range("K2").Value = "0300PM"
$ws.Range("L2").Value = "Credit"
$ws.Range("M2").Value = "cooking"
$ws.Range("N2").Value = 10
$ws.Range("O2").Value = "Active"

# --- Selection / tab state ---
$ws.Range("R5").Select()
